$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-30 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-31 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("67+13=80", $true, $false, $false, $false, $false, $true, 1, $false, "5+75=80", 2) | Out-Null
$d.Content.Find.Execute("82-51=31", $true, $false, $false, $false, $false, $true, 1, $false, "22+23=45", 2) | Out-Null
$d.Content.Find.Execute("15+1=16", $true, $false, $false, $false, $false, $true, 1, $false, "56+12=68", 2) | Out-Null
$d.Content.Find.Execute("76-33=43", $true, $false, $false, $false, $false, $true, 1, $false, "25+28=53", 2) | Out-Null
$d.Content.Find.Execute("62+6=68", $true, $false, $false, $false, $false, $true, 1, $false, "59+5=64", 2) | Out-Null
$d.Content.Find.Execute("3+88=91", $true, $false, $false, $false, $false, $true, 1, $false, "13-7=6", 2) | Out-Null
$d.Content.Find.Execute("79-40=39", $true, $false, $false, $false, $false, $true, 1, $false, "38-1=37", 2) | Out-Null
$d.Content.Find.Execute("38+35=73", $true, $false, $false, $false, $false, $true, 1, $false, "67-53=14", 2) | Out-Null
$d.Content.Find.Execute("64-31=33", $true, $false, $false, $false, $false, $true, 1, $false, "64-61=3", 2) | Out-Null
$d.Content.Find.Execute("79-9=70", $true, $false, $false, $false, $false, $true, 1, $false, "31+58=89", 2) | Out-Null
$d.Content.Find.Execute("59-28=31", $true, $false, $false, $false, $false, $true, 1, $false, "77+3=80", 2) | Out-Null
$d.Content.Find.Execute("73-60=13", $true, $false, $false, $false, $false, $true, 1, $false, "32+34=66", 2) | Out-Null
$d.Content.Find.Execute("89-10=79", $true, $false, $false, $false, $false, $true, 1, $false, "82-57=25", 2) | Out-Null
$d.Content.Find.Execute("77-49=28", $true, $false, $false, $false, $false, $true, 1, $false, "90+9=99", 2) | Out-Null
$d.Content.Find.Execute("65-27=38", $true, $false, $false, $false, $false, $true, 1, $false, "63+8=71", 2) | Out-Null
$d.Content.Find.Execute("45+8=53", $true, $false, $false, $false, $false, $true, 1, $false, "56+43=99", 2) | Out-Null
$d.Content.Find.Execute("48+33=81", $true, $false, $false, $false, $false, $true, 1, $false, "49+5=54", 2) | Out-Null
$d.Content.Find.Execute("50-28=22", $true, $false, $false, $false, $false, $true, 1, $false, "92-1=91", 2) | Out-Null
$d.Content.Find.Execute("90+4=94", $true, $false, $false, $false, $false, $true, 1, $false, "42-32=10", 2) | Out-Null
$d.Content.Find.Execute("47+7=54", $true, $false, $false, $false, $false, $true, 1, $false, "49+19=68", 2) | Out-Null
$d.Content.Find.Execute("91-91=0", $true, $false, $false, $false, $false, $true, 1, $false, "98-20=78", 2) | Out-Null
$d.Content.Find.Execute("94-13=81", $true, $false, $false, $false, $false, $true, 1, $false, "18+57=75", 2) | Out-Null
$d.Content.Find.Execute("4+45=49", $true, $false, $false, $false, $false, $true, 1, $false, "37+17=54", 2) | Out-Null
$d.Content.Find.Execute("71-53=18", $true, $false, $false, $false, $false, $true, 1, $false, "93-50=43", 2) | Out-Null
$d.Content.Find.Execute("64-39=25", $true, $false, $false, $false, $false, $true, 1, $false, "6+67=73", 2) | Out-Null
$d.Content.Find.Execute("79-5=74", $true, $false, $false, $false, $false, $true, 1, $false, "60+32=92", 2) | Out-Null
$d.Content.Find.Execute("79-39=40", $true, $false, $false, $false, $false, $true, 1, $false, "22+40=62", 2) | Out-Null
$d.Content.Find.Execute("5+13=18", $true, $false, $false, $false, $false, $true, 1, $false, "84+7=91", 2) | Out-Null
$d.Content.Find.Execute("76-16=60", $true, $false, $false, $false, $false, $true, 1, $false, "39+40=79", 2) | Out-Null
$d.Content.Find.Execute("25-16=9", $true, $false, $false, $false, $false, $true, 1, $false, "51+22=73", 2) | Out-Null
$d.Content.Find.Execute("55-34=21", $true, $false, $false, $false, $false, $true, 1, $false, "27+44=71", 2) | Out-Null
$d.Content.Find.Execute("56+34=90", $true, $false, $false, $false, $false, $true, 1, $false, "73-68=5", 2) | Out-Null
$d.Content.Find.Execute("7+40=47", $true, $false, $false, $false, $false, $true, 1, $false, "97-89=8", 2) | Out-Null
$d.Content.Find.Execute("76+12=88", $true, $false, $false, $false, $false, $true, 1, $false, "77-58=19", 2) | Out-Null
$d.Content.Find.Execute("88-28=60", $true, $false, $false, $false, $false, $true, 1, $false, "80-61=19", 2) | Out-Null
$d.Content.Find.Execute("30-4=26", $true, $false, $false, $false, $false, $true, 1, $false, "7+30=37", 2) | Out-Null
$d.Content.Find.Execute("70-3=67", $true, $false, $false, $false, $false, $true, 1, $false, "15+32=47", 2) | Out-Null
$d.Content.Find.Execute("29-8=21", $true, $false, $false, $false, $false, $true, 1, $false, "15+58=73", 2) | Out-Null
$d.Content.Find.Execute("57+22=79", $true, $false, $false, $false, $false, $true, 1, $false, "66-48=18", 2) | Out-Null
$d.Content.Find.Execute("34+43=77", $true, $false, $false, $false, $false, $true, 1, $false, "73-23=50", 2) | Out-Null
$d.Content.Find.Execute("1+91=92", $true, $false, $false, $false, $false, $true, 1, $false, "66-46=20", 2) | Out-Null
$d.Content.Find.Execute("88-1=87", $true, $false, $false, $false, $false, $true, 1, $false, "85-12=73", 2) | Out-Null
$d.Content.Find.Execute("81+14=95", $true, $false, $false, $false, $false, $true, 1, $false, "8+69=77", 2) | Out-Null
$d.Content.Find.Execute("13+7=20", $true, $false, $false, $false, $false, $true, 1, $false, "47+28=75", 2) | Out-Null
$d.Content.Find.Execute("50-40=10", $true, $false, $false, $false, $false, $true, 1, $false, "23+34=57", 2) | Out-Null
$d.Content.Find.Execute("33+20=53", $true, $false, $false, $false, $false, $true, 1, $false, "17-15=2", 2) | Out-Null
$d.Content.Find.Execute("5-3=2", $true, $false, $false, $false, $false, $true, 1, $false, "92-69=23", 2) | Out-Null
$d.Content.Find.Execute("46+12=58", $true, $false, $false, $false, $false, $true, 1, $false, "97-32=65", 2) | Out-Null
$d.Content.Find.Execute("88-37=51", $true, $false, $false, $false, $false, $true, 1, $false, "37+21=58", 2) | Out-Null
$d.Content.Find.Execute("47+43=90", $true, $false, $false, $false, $false, $true, 1, $false, "27+13=40", 2) | Out-Null
$d.Content.Find.Execute("20+5=25", $true, $false, $false, $false, $false, $true, 1, $false, "65-26=39", 2) | Out-Null
$d.Content.Find.Execute("37+28=65", $true, $false, $false, $false, $false, $true, 1, $false, "34+26=60", 2) | Out-Null
$d.Content.Find.Execute("23+53=76", $true, $false, $false, $false, $false, $true, 1, $false, "15+28=43", 2) | Out-Null
$d.Content.Find.Execute("98-21=77", $true, $false, $false, $false, $false, $true, 1, $false, "16+5=21", 2) | Out-Null
$d.Content.Find.Execute("54-47=7", $true, $false, $false, $false, $false, $true, 1, $false, "45+50=95", 2) | Out-Null
$d.Content.Find.Execute("43-41=2", $true, $false, $false, $false, $false, $true, 1, $false, "89-0=89", 2) | Out-Null
$d.Content.Find.Execute("4+2=6", $true, $false, $false, $false, $false, $true, 1, $false, "40+3=43", 2) | Out-Null
$d.Content.Find.Execute("33-6=27", $true, $false, $false, $false, $false, $true, 1, $false, "96+1=97", 2) | Out-Null
$d.Content.Find.Execute("50-34=16", $true, $false, $false, $false, $false, $true, 1, $false, "88-71=17", 2) | Out-Null
$d.Content.Find.Execute("65+33=98", $true, $false, $false, $false, $false, $true, 1, $false, "70-24=46", 2) | Out-Null
$d.Content.Find.Execute("68+0=68", $true, $false, $false, $false, $false, $true, 1, $false, "52-24=28", 2) | Out-Null
$d.Content.Find.Execute("36-7=29", $true, $false, $false, $false, $false, $true, 1, $false, "28+64=92", 2) | Out-Null
$d.Content.Find.Execute("51+11=62", $true, $false, $false, $false, $false, $true, 1, $false, "31+17=48", 2) | Out-Null
$d.Content.Find.Execute("60+13=73", $true, $false, $false, $false, $false, $true, 1, $false, "78-29=49", 2) | Out-Null
$d.Content.Find.Execute("78-2=76", $true, $false, $false, $false, $false, $true, 1, $false, "7+7=14", 2) | Out-Null
$d.Content.Find.Execute("0+15=15", $true, $false, $false, $false, $false, $true, 1, $false, "90-60=30", 2) | Out-Null
$d.Content.Find.Execute("98-93=5", $true, $false, $false, $false, $false, $true, 1, $false, "19+43=62", 2) | Out-Null
$d.Content.Find.Execute("88-27=61", $true, $false, $false, $false, $false, $true, 1, $false, "90-56=34", 2) | Out-Null
$d.Content.Find.Execute("13+56=69", $true, $false, $false, $false, $false, $true, 1, $false, "2-0=2", 2) | Out-Null
$d.Content.Find.Execute("51-3=48", $true, $false, $false, $false, $false, $true, 1, $false, "4+49=53", 2) | Out-Null
$d.Content.Find.Execute("92-20=72", $true, $false, $false, $false, $false, $true, 1, $false, "69-19=50", 2) | Out-Null
$d.Content.Find.Execute("68-19=49", $true, $false, $false, $false, $false, $true, 1, $false, "7+81=88", 2) | Out-Null
$d.Content.Find.Execute("30+52=82", $true, $false, $false, $false, $false, $true, 1, $false, "71-5=66", 2) | Out-Null
$d.Content.Find.Execute("48+41=89", $true, $false, $false, $false, $false, $true, 1, $false, "24+12=36", 2) | Out-Null
$d.Content.Find.Execute("17+36=53", $true, $false, $false, $false, $false, $true, 1, $false, "23+17=40", 2) | Out-Null
$d.Content.Find.Execute("90-39=51", $true, $false, $false, $false, $false, $true, 1, $false, "85+0=85", 2) | Out-Null
$d.Content.Find.Execute("22+65=87", $true, $false, $false, $false, $false, $true, 1, $false, "37-27=10", 2) | Out-Null
$d.Content.Find.Execute("42-5=37", $true, $false, $false, $false, $false, $true, 1, $false, "55-18=37", 2) | Out-Null
$d.Content.Find.Execute("71+28=99", $true, $false, $false, $false, $false, $true, 1, $false, "91-14=77", 2) | Out-Null
$d.Content.Find.Execute("58+18=76", $true, $false, $false, $false, $false, $true, 1, $false, "37+15=52", 2) | Out-Null
$d.Content.Find.Execute("53-14=39", $true, $false, $false, $false, $false, $true, 1, $false, "53-44=9", 2) | Out-Null
$d.Content.Find.Execute("99-32=67", $true, $false, $false, $false, $false, $true, 1, $false, "57-11=46", 2) | Out-Null
$d.Content.Find.Execute("41-14=27", $true, $false, $false, $false, $false, $true, 1, $false, "33+60=93", 2) | Out-Null
$d.Content.Find.Execute("76-42=34", $true, $false, $false, $false, $false, $true, 1, $false, "8-2=6", 2) | Out-Null
$d.Content.Find.Execute("53+29=82", $true, $false, $false, $false, $false, $true, 1, $false, "13+33=46", 2) | Out-Null
$d.Content.Find.Execute("18+29=47", $true, $false, $false, $false, $false, $true, 1, $false, "42+31=73", 2) | Out-Null
$d.Content.Find.Execute("85-69=16", $true, $false, $false, $false, $false, $true, 1, $false, "15+80=95", 2) | Out-Null
$d.Content.Find.Execute("18+53=71", $true, $false, $false, $false, $false, $true, 1, $false, "49+6=55", 2) | Out-Null
$d.Content.Find.Execute("83+8=91", $true, $false, $false, $false, $false, $true, 1, $false, "42+42=84", 2) | Out-Null
$d.Content.Find.Execute("21+7=28", $true, $false, $false, $false, $false, $true, 1, $false, "71-0=71", 2) | Out-Null
$d.Content.Find.Execute("85-42=43", $true, $false, $false, $false, $false, $true, 1, $false, "57-56=1", 2) | Out-Null
$d.Content.Find.Execute("45-33=12", $true, $false, $false, $false, $false, $true, 1, $false, "24+34=58", 2) | Out-Null
$d.Content.Find.Execute("48+38=86", $true, $false, $false, $false, $false, $true, 1, $false, "34-9=25", 2) | Out-Null
$d.Content.Find.Execute("1+41=42", $true, $false, $false, $false, $false, $true, 1, $false, "65-5=60", 2) | Out-Null
$d.Content.Find.Execute("58+26=84", $true, $false, $false, $false, $false, $true, 1, $false, "92-24=68", 2) | Out-Null
$d.Content.Find.Execute("99-12=87", $true, $false, $false, $false, $false, $true, 1, $false, "17+61=78", 2) | Out-Null
$d.Content.Find.Execute("11-11=0", $true, $false, $false, $false, $false, $true, 1, $false, "32-6=26", 2) | Out-Null
$d.Content.Find.Execute("67-19=48", $true, $false, $false, $false, $false, $true, 1, $false, "57+8=65", 2) | Out-Null
$d.Content.Find.Execute("25+53=78", $true, $false, $false, $false, $false, $true, 1, $false, "45+40=85", 2) | Out-Null
$d.Content.Find.Execute("51+44=95", $true, $false, $false, $false, $false, $true, 1, $false, "88-6=82", 2) | Out-Null
